# Fix formula bugs when the formula isn't the first text in the cell.
# On the "MultiLevel" sheet, update a few cells so that the jt:if tag
# markup wraps around the formula text rather than the formula being
# the whole/only content, and fully-qualify the two COUNTA formula
# references with absolute cell references.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultiLevel")

$ws.Range("B5").Value = '<jt:if test="true">$[COUNTA(B3||$Z$1)]'
$ws.Range("F5").Value = '$[SUM(C3)/SUM(E3||1)]</jt:if>'
$ws.Range("I1").Value = '$[COUNTA(''Formula Test''!$E$3)]'
$ws.Range("I2").Value = '$[COUNTA(''Formula Test''!$K$3)]'
